# Run30 Channel_Map / Run_Parameters update:
# Add a new "simulations" sheet (a copy of the last existing run-parameters
# sheet, "20200924_Afternoon_AfterSeventh") with an updated Baseline Length
# value, and move the active-tab/selection state onto the new sheet.

$wb = $excel.ActiveWorkbook

# Last existing sheet ("20200924_Afternoon_AfterSeventh") is the template
# for the new simulation-parameters sheet.
$src = $wb.Worksheets.Item($wb.Worksheets.Count)

# Copy it to just after itself (i.e. append at the end of the tab strip).
$src.Copy([System.Reflection.Missing]::Value, $src) | Out-Null

# The freshly inserted copy becomes the active sheet.
$new = $wb.ActiveSheet
$new.Name = "simulations"

# Simulation run uses a shorter baseline length (row 10, "Baseline Length [samples]").
$new.Range("B10").Value = 125

# Leave the selection on the old parameters sheet at B10 ...
$src.Range("B10").Select() | Out-Null

# ... and select B9 on the new "simulations" sheet, which stays the active tab.
$new.Activate()
$new.Range("B9").Select() | Out-Null
